$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.823.10"
$ws.Range("E2").Value = "  +4.80%  "
$ws.Range("D3").Value = "2.271.57"
$ws.Range("E3").Value = "  +1.91%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.627"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.25"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.56%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.420"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.95"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0934"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.96%  "
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").Value = "2.604.37"
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.92%  "
$ws.Range("E16").Value = "  +3.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.808"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").Value = "2.261.34"
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("D19").Value = "43.745.48"
$ws.Range("E19").Value = "  +4.90%  "
$ws.Range("D20").Value = "0.0₃0936"
$ws.Range("E20").Value = "  +3.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.25%  "
$ws.Range("E22").Value = "  +2.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.33%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "170.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("E29").Value = "  -1.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.44"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.38%  "
$ws.Range("E32").Value = "  -0.59%  "
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.06"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.77"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.94%  "
$ws.Range("E36").Value = "  +3.02%  "
$ws.Range("E37").Value = "  -2.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.87%  "
$ws.Range("E40").Value = "  +4.30%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.72"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.000226"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -10.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0986"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.55%  "
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "98.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("D48").Value = "1.476.08"
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.67"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.36%  "
